$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F ("dSF")
$changes = @{
    2  = -3
    3  = 2
    7  = 9
    8  = 0
    11 = -4
    15 = -3
    17 = -3
    19 = 1
    22 = -2
    23 = -6
    26 = 3
    27 = 7
    33 = 13
    38 = 1
    45 = -1
    47 = -3
    49 = -4
    51 = -3
    52 = -2
    55 = -1
    60 = 1
    63 = 2
    65 = 0
    67 = -3
    70 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
